# issue #5: add legislator_id, name, date into dataframe
#
# This adds three new columns (date, legislator_name, legislator_id) to the
# "股票" (stocks) worksheet, and leaves the "保險" (insurance) worksheet data
# values unchanged (its column layout already had owner info in column D).

$wb = $excel.ActiveWorkbook

# --- Sheet "股票" (stocks) : add date / legislator_name / legislator_id columns ---
$wsStock = $wb.Worksheets.Item("股票")

# Header row (row 1): new headers in H, I, J
$wsStock.Range("H1").Value = "date"
$wsStock.Range("I1").Value = "legislator_name"
$wsStock.Range("J1").Value = "legislator_id"

# Copy the header style (bold / border / centered) used by the other header cells
$wsStock.Range("G1").Copy() | Out-Null
$wsStock.Range("H1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Data rows (2 through 13): fill in the date, legislator name and id.
# Force column H to text format first so the date-like string "2012-04-24"
# is kept as text instead of being auto-converted into a date serial number.
$lastRow = 13
$wsStock.Range("H2:H$lastRow").NumberFormat = "@"
for ($r = 2; $r -le $lastRow; $r++) {
    $wsStock.Range("H$r").Value = "2012-04-24"
    $wsStock.Range("I$r").Value = "段宜康"
    $wsStock.Range("J$r").Value = 917
}

$excel.CutCopyMode = 0
